# The "waitTime::6000" value was a mistake: the wait time delay is
# interpreted in seconds, not milliseconds, so replace every occurrence
# of that value (in the t-omdb&imdb sheet) with the corrected "waitTime::6".
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("t-omdb&imdb")

$ws.Cells.Replace("waitTime::6000", "waitTime::6")

# Re-select the cell that was left selected in the sheet after the edit.
$ws.Activate()
$ws.Range("F7").Select()
